$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.154.36"
$ws.Range("E2").Value = "  +0.37%  "

# Row 3
$ws.Range("D3").Value = "3.340.29"
$ws.Range("E3").Value = "  +0.91%  "

# Row 4
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.12%  "

# Row 5
$ws.Range("D5").Value = "'585.88"
$ws.Range("E5").Value = "  +5.42%  "

# Row 6
$ws.Range("D6").Value = "'185.39"
$ws.Range("E6").Value = "  -1.29%  "

# Row 7
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.04%  "

# Row 8
$ws.Range("E8").Value = "  -1.44%  "

# Row 9
$ws.Range("E9").Value = "  -1.01%  "

# Row 10
$ws.Range("E10").Value = "  -0.77%  "

# Row 11
$ws.Range("D11").Value = "'46.89"
$ws.Range("E11").Value = "  -0.82%  "

# Row 12
$ws.Range("D12").Value = "'0.0000269"
$ws.Range("E12").Value = "  -0.33%  "

# Row 13
$ws.Range("D13").Value = "'655.71"
$ws.Range("E13").Value = "  +8.31%  "

# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'8.48"
$ws.Range("E14").Value = "  -2.54%  "

# Row 15
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "3.629.29"
$ws.Range("E15").Value = "  -5.32%  "

# Row 16
$ws.Range("D16").Value = "66.347.56"
$ws.Range("E16").Value = "  +0.72%  "

# Row 17
$ws.Range("E17").Value = "  -0.21%  "

# Row 18
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").Value = "'17.86"
$ws.Range("E18").Value = "  -0.75%  "

# Row 19
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "3.334.86"
$ws.Range("E19").Value = "  +0.92%  "

# Row 20
$ws.Range("D20").Value = "'11.10"
$ws.Range("E20").Value = "  +0.17%  "

# Row 21
$ws.Range("D21").Value = "'0.897"
$ws.Range("E21").Value = "  -1.17%  "

# Row 22
$ws.Range("E22").Value = "  -4.40%  "

# Row 23
$ws.Range("E23").Value = "  -0.55%  "

# Row 24
$ws.Range("D24").Value = "'100.04"
$ws.Range("E24").Value = "  +0.05%  "

# Row 25
$ws.Range("D25").Value = "'4.00"
$ws.Range("E25").Value = "  +1.24%  "

# Row 26
$ws.Range("E26").Value = "  +0.84%  "

# Row 27
$ws.Range("D27").Value = "'9.59"
$ws.Range("E27").Value = "  +0.12%  "

# Row 28
$ws.Range("D28").Value = "'32.06"
$ws.Range("E28").Value = "  +5.71%  "

# Row 29
$ws.Range("D29").Value = "'8.52"
$ws.Range("E29").Value = "  -1.92%  "

# Row 30
$ws.Range("D30").Value = "'6.75"
$ws.Range("E30").Value = "  +0.32%  "

# Row 31
$ws.Range("D31").Value = "'604.84"
$ws.Range("E31").Value = "  +3.83%  "

# Row 32
$ws.Range("D32").Value = "'3.86"
$ws.Range("E32").Value = "  +0.29%  "

# Row 33
$ws.Range("D33").Value = "'11.08"
$ws.Range("E33").Value = "  -0.22%  "

# Row 34
$ws.Range("D34").Value = "3.881.89"
$ws.Range("E34").Value = "  +4.87%  "

# Row 35
$ws.Range("E35").Value = "  +0.40%  "

# Row 36
$ws.Range("E36").Value = "  +0.00%  "

# Row 37
$ws.Range("D37").Value = "'56.39"
$ws.Range("E37").Value = "  -1.10%  "

# Row 38
$ws.Range("D38").Value = "'2.77"
$ws.Range("E38").Value = "  +3.40%  "

# Row 39
$ws.Range("D39").Value = "'0.130"
$ws.Range("E39").Value = "  -0.43%  "

# Row 40
$ws.Range("B40").Value = "PEPE"
$ws.Range("C40").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D40").Value = "0.0₃0699"
$ws.Range("E40").Value = "  -3.23%  "

# Row 41
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").Value = "'32.90"
$ws.Range("E41").Value = "  -3.20%  "

# Row 42
$ws.Range("E42").Value = "  -2.66%  "

# Row 43
$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D43").Value = "'0.340"
$ws.Range("E43").Value = "  -0.11%  "

# Row 44
$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").Value = "'3.37"
$ws.Range("E44").Value = "  +0.20%  "

# Row 45
$ws.Range("E45").Value = "  -1.03%  "

# Row 46
$ws.Range("E46").Value = "  -1.25%  "

# Row 47
$ws.Range("E47").Value = "  +0.38%  "

# Row 48
$ws.Range("E48").Value = "  -1.57%  "

# Row 49
$ws.Range("D49").Value = "'2.84"
$ws.Range("E49").Value = "  -18.12%  "

# Row 50
$ws.Range("E50").Value = "  +5.93%  "

# Row 51
$ws.Range("D51").Value = "'129.39"
$ws.Range("E51").Value = "  +4.59%  "

# Reset style on text-forced numeric-looking cells to avoid quotePrefix style residue
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"
